# Apply the TC_101 "loc to dev" edits to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the sheet: "My Series" -> "Data"
# ---------------------------------------------------------------------
$ws.Name = "Data"

# ---------------------------------------------------------------------
# 2) Header row (row 1): reorder the series names and tag each with
#    "[ACCUMULATE()]".
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Govt Revenue - Tax ; Individual Income [ACCUMULATE()]"
$ws.Range("C1").Value = "Value Added in Industry [ACCUMULATE()]"
$ws.Range("D1").Value = "Govt Revenue [ACCUMULATE()]"

# ---------------------------------------------------------------------
# 3) Row 5 (Unit): swap C5 / D5.
# ---------------------------------------------------------------------
$ws.Range("C5").Value = "LKR mn"
$ws.Range("D5").Value = "RMB mn"

# ---------------------------------------------------------------------
# 4) Row 8 (Series ID): swap C8 / D8.
# ---------------------------------------------------------------------
$ws.Range("C8").Value = 310902601
$ws.Range("D8").Value = 310901801

# ---------------------------------------------------------------------
# 5) Row 9 (SR Code): swap C9 / D9.
# ---------------------------------------------------------------------
$ws.Range("C9").Value = "SR4825076"
$ws.Range("D9").Value = "SR4825071"

# ---------------------------------------------------------------------
# 6) Row 11 label: "Function Description" -> "Function Information".
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "Function Information"

# ---------------------------------------------------------------------
# 7) Row 12 (First Obs. Date): swap C12 / D12.
# ---------------------------------------------------------------------
$ws.Range("C12").Value = 31017
$ws.Range("D12").Value = 18598

# ---------------------------------------------------------------------
# 8) Row 13 (Last Obs. Date): swap C13 / D13.
# ---------------------------------------------------------------------
$ws.Range("C13").Value = 39052
$ws.Range("D13").Value = 41244

# ---------------------------------------------------------------------
# 9) Rows 17-25 (summary statistics): swap C / D. A couple of values
#    were refreshed (1-ULP float jitter from the resummarised series)
#    so the exact target literals from the diff are used below.
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 1453915.3
$ws.Range("D17").Value = 37803803.9

$ws.Range("C18").Value = 407944167687.1222
$ws.Range("D18").Value = 310352334867253.9

$ws.Range("C19").Value = 638705.0709733892
$ws.Range("D19").Value = 17616819.65813506

$ws.Range("C20").Value = 0.4485025925517874
$ws.Range("D20").Value = 0.421627230980701

$ws.Range("B21").Value = -1.129868563809901
$ws.Range("C21").Value = -0.9169323140962309
$ws.Range("D21").Value = -1.334739153466928

$ws.Range("C22").Value = 0.4393000548060738
$ws.Range("D22").Value = 0.4660065348115686

$ws.Range("C23").Value = 643058
$ws.Range("D23").Value = 16641896

$ws.Range("C24").Value = 2549974
$ws.Range("D24").Value = 63990353

$ws.Range("C25").Value = 1362529.5
$ws.Range("D25").Value = 34521187.5

# ---------------------------------------------------------------------
# 10) Rows 27-42 (raw observations): the whole C/D block is swapped,
#     including blank cells.
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 643058
$ws.Range("D27").ClearContents()

$ws.Range("C28").Value = 772535
$ws.Range("D28").ClearContents()

$ws.Range("C29").Value = 915947
$ws.Range("D29").ClearContents()

$ws.Range("C30").Value = 1083440
$ws.Range("D30").ClearContents()

$ws.Range("C31").Value = 1262998
$ws.Range("D31").ClearContents()

$ws.Range("C32").Value = 1462061
$ws.Range("D32").ClearContents()

$ws.Range("C33").Value = 1683077
$ws.Range("D33").Value = 16641896

$ws.Range("C34").Value = 1939649
$ws.Range("D34").Value = 19281543

$ws.Range("C35").Value = 2226414
$ws.Range("D35").Value = 22446472

$ws.Range("C36").Value = 2549974
$ws.Range("D36").Value = 26322492

$ws.Range("C37").ClearContents()
$ws.Range("D37").Value = 31454670

$ws.Range("C38").ClearContents()
$ws.Range("D38").Value = 37587705

$ws.Range("C39").ClearContents()
$ws.Range("D39").Value = 44439535

$ws.Range("C40").ClearContents()
$ws.Range("D40").Value = 52749686

$ws.Range("C41").ClearContents()
$ws.Range("D41").Value = 63123687

$ws.Range("C42").ClearContents()
$ws.Range("D42").Value = 63990353

# ---------------------------------------------------------------------
# 11) Number format for the observation block: "0.000" -> "###0.000".
#     All 48 cells in B27:D42 share this format, so re-apply it across
#     the whole block in one shot.
# ---------------------------------------------------------------------
$ws.Range("B27:D42").NumberFormat = "###0.000"
